# Fix error in high current MX fab BOM: the resistor values in column A (Comment)
# and column D (LCSC Part #) for designators R14, R15, R16 had been entered one row
# off - row 16 (R14) should show 30k instead of 300k, row 17 (R15) should show 3k
# instead of 30k, row 18 (R16) should show 300 Ohm instead of 3k, and row 19 (R17)
# needs a brand-new "30.1 Ohm Resistor, 1%, 200mW" / " C365143" part (Qty cell left
# blank because this new part hasn't been quoted yet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 -> R14: now 30k Resistor, 1% / C22984
$ws.Range("A16").Value = "30k Resistor, 1%"
$ws.Range("D16").Value = "C22984"

# Row 17 -> R15: now 3k Resistor, 1% / C4211
$ws.Range("A17").Value = "3k Resistor, 1%"
$ws.Range("D17").Value = "C4211"

# Row 18 -> R16: now 300 Ohm Resistor, 1% / C23025
$ws.Range("A18").Value = "300 Ohm Resistor, 1%"
$ws.Range("D18").Value = "C23025"

# Row 19 -> R17: now 30.1 Ohm Resistor, 1%, 200mW / " C365143", Qty cell removed
$ws.Range("A19").Value = "30.1 Ohm Resistor, 1%, 200mW"
$ws.Range("D19").Value = " C365143"
$ws.Range("E19").Clear()

# Update the remembered selection to match the author's saved cursor position
$ws.Range("A20").Select()
